$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 325.5
$ws.Range("I33").Value = 280.53333
$ws.Range("K33").Value = 280.53333
$ws.Range("M33").Value = -51.53332999999998
$ws.Range("H61").Value = 2404.6667
$ws.Range("I61").Value = 2404.6667
$ws.Range("K61").Value = 7214.000100000001
$ws.Range("M61").Value = -7042.000100000001
$ws.Range("H132").Value = 1173.4746
$ws.Range("I132").Value = 1173.9122
$ws.Range("K132").Value = 3521.7366
$ws.Range("M132").Value = -991.7366000000002
$ws.Range("H137").Value = 27780900
$ws.Range("J137").Value = 3735.8948
$ws.Range("L137").Value = 11207.6844
$ws.Range("N137").Value = -16307.6844
$ws.Range("H138").Value = 3298.2354
$ws.Range("I138").Value = 1969.68
$ws.Range("J138").Value = 6988.6665
$ws.Range("K138").Value = 5909.04
$ws.Range("L138").Value = 20965.9995
$ws.Range("M138").Value = -769.04
$ws.Range("N138").Value = -31245.9995
$ws.Range("H141").Value = 788.6667
$ws.Range("I141").Value = 788.6667
$ws.Range("K141").Value = 2366.0001
$ws.Range("M141").Value = 2813.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7711.7856
$ws.Range("I2").Value = 612.6923
$ws.Range("K2").Value = 612.6923
$ws.Range("M2").Value = -499.6923
$ws.Range("H32").Value = 10640501
$ws.Range("I32").Value = 10990869
$ws.Range("K32").Value = 10990869
$ws.Range("M32").Value = -10990582
$ws.Range("H61").Value = 4998.595
$ws.Range("I61").Value = 5184.8
$ws.Range("J61").Value = 1274.5
$ws.Range("K61").Value = 5184.8
$ws.Range("L61").Value = 1274.5
$ws.Range("M61").Value = -4972.8
$ws.Range("N61").Value = -1698.5
$ws.Range("H116").Value = 7711.7856
$ws.Range("I116").Value = 612.6923
$ws.Range("K116").Value = 612.6923
$ws.Range("M116").Value = 1681.3077
$ws.Range("H132").Value = 7070.643
$ws.Range("I132").Value = 2854.2778
$ws.Range("K132").Value = 8562.8334
$ws.Range("M132").Value = -6032.8334
$ws.Range("H136").Value = 4998.595
$ws.Range("I136").Value = 5184.8
$ws.Range("J136").Value = 1274.5
$ws.Range("K136").Value = 15554.4
$ws.Range("L136").Value = 3823.5
$ws.Range("M136").Value = -13004.4
$ws.Range("N136").Value = -8923.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7711.7856
$ws.Range("I3").Value = 612.6923
$ws.Range("K3").Value = 612.6923
$ws.Range("M3").Value = -498.6923
$ws.Range("H94").Value = 1183.7778
$ws.Range("J94").Value = 1248.25
$ws.Range("L94").Value = 1248.25
$ws.Range("N94").Value = -2150.25
$ws.Range("H105").Value = 17727.5
$ws.Range("I105").Value = 16578.285
$ws.Range("K105").Value = 16578.285
$ws.Range("M105").Value = -14831.285
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2195.05
$ws.Range("I16").Value = 1438.2222
$ws.Range("J16").Value = 9006.5
$ws.Range("K16").Value = 1438.2222
$ws.Range("L16").Value = 9006.5
$ws.Range("M16").Value = -1151.2222
$ws.Range("N16").Value = -9580.5
$ws.Range("H22").Value = 1882.0834
$ws.Range("I22").Value = 308.33334
$ws.Range("K22").Value = 308.33334
$ws.Range("M22").Value = 41.66665999999998
$ws.Range("H94").Value = 2582.0908
$ws.Range("J94").Value = 3285.4285
$ws.Range("L94").Value = 3285.4285
$ws.Range("N94").Value = -4187.4285
$ws.Range("H107").Value = 704.05884
$ws.Range("I107").Value = 508
$ws.Range("K107").Value = 508
$ws.Range("M107").Value = 1412
$ws.Range("H113").Value = 2195.05
$ws.Range("I113").Value = 1438.2222
$ws.Range("J113").Value = 9006.5
$ws.Range("K113").Value = 1438.2222
$ws.Range("L113").Value = 9006.5
$ws.Range("M113").Value = 731.7778000000001
$ws.Range("N113").Value = -13346.5
$ws.Range("H122").Value = 6664.7827
$ws.Range("I122").Value = 2769.2354
$ws.Range("J122").Value = 17702.166
$ws.Range("K122").Value = 8307.706200000001
$ws.Range("L122").Value = 53106.49800000001
$ws.Range("M122").Value = -5857.706200000001
$ws.Range("N122").Value = -58006.49800000001
$ws.Range("H132").Value = 3590.762
$ws.Range("I132").Value = 1807.7142
$ws.Range("K132").Value = 5423.142599999999
$ws.Range("M132").Value = -2893.142599999999
$ws.Range("H134").Value = 4490.5654
$ws.Range("I134").Value = 2255.543
$ws.Range("K134").Value = 6766.629000000001
$ws.Range("M134").Value = -4231.629000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 934.5
$ws.Range("I59").Value = 934.5
$ws.Range("K59").Value = 2803.5
$ws.Range("M59").Value = -2263.5
$ws.Range("H74").Value = 12007.5
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -13939
$ws.Range("H77").Value = 12007.5
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -39696
$ws.Range("H132").Value = 3813.7666
$ws.Range("I132").Value = 3261.3333
$ws.Range("J132").Value = 4182.0557
$ws.Range("K132").Value = 29351.9997
$ws.Range("L132").Value = 37638.5013
$ws.Range("M132").Value = -26821.9997
$ws.Range("N132").Value = -42698.5013
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9505.24
$ws.Range("I70").Value = 7213.3125
$ws.Range("J70").Value = 13579.777
$ws.Range("K70").Value = 7213.3125
$ws.Range("L70").Value = 13579.777
$ws.Range("M70").Value = -6943.3125
$ws.Range("N70").Value = -14119.777
$ws.Range("H73").Value = 9505.24
$ws.Range("I73").Value = 7213.3125
$ws.Range("J73").Value = 13579.777
$ws.Range("K73").Value = 7213.3125
$ws.Range("L73").Value = 13579.777
$ws.Range("M73").Value = -6277.3125
$ws.Range("N73").Value = -15451.777
$ws.Range("H107").Value = 1256.4
$ws.Range("I107").Value = 337.7143
$ws.Range("J107").Value = 3400
$ws.Range("K107").Value = 337.7143
$ws.Range("L107").Value = 3400
$ws.Range("M107").Value = 1582.2857
$ws.Range("N107").Value = -7240
$ws.Range("H113").Value = 1087.875
$ws.Range("I113").Value = 1087.875
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1087.875
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1082.125
$ws.Range("N113").ClearContents()
$ws.Range("H130").Value = 72245
$ws.Range("J130").Value = 72245
$ws.Range("L130").Value = 72245
$ws.Range("N130").Value = -82285
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 302.36
$ws.Range("I16").Value = 289.95834
$ws.Range("K16").Value = 289.95834
$ws.Range("M16").Value = -119.95834
$ws.Range("H22").Value = 19001
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 19001
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H109").Value = 46052.5
$ws.Range("J109").Value = 46052.5
$ws.Range("L109").Value = 46052.5
$ws.Range("N109").Value = -48826.5
$ws.Range("H122").Value = 5410.4375
$ws.Range("J122").Value = 7745.125
$ws.Range("L122").Value = 23235.375
$ws.Range("N122").Value = -28135.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1132.7142
$ws.Range("I107").Value = 934.25
$ws.Range("K107").Value = 2802.75
$ws.Range("M107").Value = -882.75
$ws.Range("H122").Value = 3162.25
$ws.Range("I122").Value = 1425.9412
$ws.Range("K122").Value = 4277.8236
$ws.Range("M122").Value = -1827.8236
$ws.Range("H132").Value = 4103.175
$ws.Range("I132").Value = 1677.3043
$ws.Range("J132").Value = 10667.294
$ws.Range("K132").Value = 5031.9129
$ws.Range("L132").Value = 32001.882
$ws.Range("M132").Value = -2501.9129
$ws.Range("N132").Value = -37061.882

Write-Output "Applied all changes"